# Insert a new weekly record at row 472, pushing the existing rows
# (old 472..496) down to (473..497). This matches the diff: dimension
# grows from A1:R496 to A1:R497, and a brand-new "Ají" record appears at
# row 472 while everything that used to be at row N now lives at row N+1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 472:496 down to 473:497 by inserting a blank row at 472.
$ws.Rows.Item(472).Insert()

# Populate the newly inserted row 472 with the new record's data.
$ws.Cells.Item(472, 1).Value = 3
$ws.Cells.Item(472, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(472, 3).Value = "Coquimbo"
$ws.Cells.Item(472, 4).Value = 44753
$ws.Cells.Item(472, 5).Value = 5
$ws.Cells.Item(472, 6).Value = 100112021
$ws.Cells.Item(472, 7).Value = "Ají"
$ws.Cells.Item(472, 8).Value = "Inferno"
$ws.Cells.Item(472, 9).Value = "Primera"
$ws.Cells.Item(472, 10).Value = 73
$ws.Cells.Item(472, 11).Value = 15000
$ws.Cells.Item(472, 12).Value = 16000
$ws.Cells.Item(472, 13).Value = 15479
$ws.Cells.Item(472, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(472, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(472, 16).Value = 1032
$ws.Cells.Item(472, 17).Value = 15
$ws.Cells.Item(472, 18).Value = "Hortaliza"
